$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7817.407
$ws.Range("I74").Value = 6098.2666
$ws.Range("K74").Value = 6098.2666
$ws.Range("M74").Value = -5162.2666
$ws.Range("H77").Value = 7817.407
$ws.Range("I77").Value = 6098.2666
$ws.Range("K77").Value = 30491.333
$ws.Range("M77").Value = -25811.333
$ws.Range("H98").Value = 1629.1892
$ws.Range("I98").Value = 1723.1515
$ws.Range("J98").Value = 854
$ws.Range("K98").Value = 1723.1515
$ws.Range("L98").Value = 854
$ws.Range("M98").Value = -225.1514999999999
$ws.Range("N98").Value = -3850
$ws.Range("H99").Value = 1854.3334
$ws.Range("J99").Value = 1898.75
$ws.Range("L99").Value = 5696.25
$ws.Range("N99").Value = -8692.25
$ws.Range("H101").Value = 2996.3333
$ws.Range("J101").Value = 2795
$ws.Range("L101").Value = 8385
$ws.Range("N101").Value = -11629
$ws.Range("H122").Value = 1629.1892
$ws.Range("I122").Value = 1723.1515
$ws.Range("J122").Value = 854
$ws.Range("K122").Value = 5169.4545
$ws.Range("L122").Value = 2562
$ws.Range("M122").Value = -2719.4545
$ws.Range("N122").Value = -7462
$ws.Range("H135").Value = 1491.4
$ws.Range("I135").Value = 1533.7391
$ws.Range("J135").Value = 1352.2858
$ws.Range("K135").Value = 13803.6519
$ws.Range("L135").Value = 12170.5722
$ws.Range("M135").Value = -11268.6519
$ws.Range("N135").Value = -17240.5722
$ws.Range("H137").Value = 4040.1738
$ws.Range("J137").Value = 4304.7896
$ws.Range("L137").Value = 12914.3688
$ws.Range("N137").Value = -18014.3688
$ws.Range("H141").Value = 1335.5714
$ws.Range("I141").Value = 1335.5714
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4006.7142
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1173.2858
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5058862
$ws.Range("I32").Value = 6066415.5
$ws.Range("K32").Value = 6066415.5
$ws.Range("M32").Value = -6066128.5
$ws.Range("H74").Value = 3475180.8
$ws.Range("I74").Value = 4311920.5
$ws.Range("J74").Value = 8687.143
$ws.Range("K74").Value = 4311920.5
$ws.Range("L74").Value = 8687.143
$ws.Range("M74").Value = -4311046.5
$ws.Range("N74").Value = -10435.143
$ws.Range("H77").Value = 3475180.8
$ws.Range("I77").Value = 4311920.5
$ws.Range("J77").Value = 8687.143
$ws.Range("K77").Value = 21559602.5
$ws.Range("L77").Value = 43435.715
$ws.Range("M77").Value = -21555234.5
$ws.Range("N77").Value = -52171.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 44999.75
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 44999.75
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 44999.75
$ws.Range("M59").Value = $null
$ws.Range("N59").Value = -47289.75
$ws.Range("H74").Value = 44894.42
$ws.Range("I74").Value = 44494
$ws.Range("K74").Value = 44494
$ws.Range("M74").Value = -43620
$ws.Range("H77").Value = 44894.42
$ws.Range("I77").Value = 44494
$ws.Range("K77").Value = 133482
$ws.Range("M77").Value = -129114
$ws.Range("H94").Value = 19884.908
$ws.Range("J94").Value = 2252.1667
$ws.Range("L94").Value = 2252.1667
$ws.Range("N94").Value = -3154.1667
$ws.Range("H132").Value = 37862.727
$ws.Range("I132").Value = 89139.25
$ws.Range("J132").Value = 8561.857
$ws.Range("K132").Value = 267417.75
$ws.Range("L132").Value = 25685.571
$ws.Range("M132").Value = -264887.75
$ws.Range("N132").Value = -30745.571
$ws.Range("H134").Value = 3599.2886
$ws.Range("I134").Value = 1303.26
$ws.Range("K134").Value = 3909.78
$ws.Range("M134").Value = -1374.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I68").Value = 2747.25
$ws.Range("J68").Value = 10381.909
$ws.Range("K68").Value = 8241.75
$ws.Range("L68").Value = 31145.727
$ws.Range("M68").Value = -7430.75
$ws.Range("N68").Value = -32767.727
$ws.Range("I71").Value = 2747.25
$ws.Range("J71").Value = 10381.909
$ws.Range("K71").Value = 24725.25
$ws.Range("L71").Value = 93437.181
$ws.Range("M71").Value = -20669.25
$ws.Range("N71").Value = -101549.181
$ws.Range("H76").Value = 7013
$ws.Range("I76").Value = 7013
$ws.Range("K76").Value = 21039
$ws.Range("M76").Value = -20656
$ws.Range("H79").Value = 7013
$ws.Range("I79").Value = 7013
$ws.Range("K79").Value = 21039
$ws.Range("M79").Value = -19713
$ws.Range("H81").Value = 5507.2
$ws.Range("J81").Value = 6444.875
$ws.Range("L81").Value = 19334.625
$ws.Range("N81").Value = -21580.625
$ws.Range("H84").Value = 5507.2
$ws.Range("J84").Value = 6444.875
$ws.Range("L84").Value = 58003.875
$ws.Range("N84").Value = -69235.875
$ws.Range("H122").Value = 776.0540999999999
$ws.Range("J122").Value = 830.71875
$ws.Range("L122").Value = 7476.46875
$ws.Range("N122").Value = -12376.46875
$ws.Range("H132").Value = 3821.875
$ws.Range("I132").Value = 2997.8
$ws.Range("J132").Value = 4196.4546
$ws.Range("K132").Value = 26980.2
$ws.Range("L132").Value = 37768.0914
$ws.Range("M132").Value = -24450.2
$ws.Range("N132").Value = -42828.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1661.3636
$ws.Range("I132").Value = 1726.1
$ws.Range("K132").Value = 5178.299999999999
$ws.Range("M132").Value = -2648.299999999999
$ws.Range("H135").Value = 97496
$ws.Range("J135").Value = 97496
$ws.Range("L135").Value = 97496
$ws.Range("N135").Value = -107636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3330
$ws.Range("I40").Value = 2995
$ws.Range("K40").Value = 2995
$ws.Range("M40").Value = -2859
$ws.Range("H46").Value = 2965.7932
$ws.Range("I46").Value = 1228.6
$ws.Range("K46").Value = 1228.6
$ws.Range("M46").Value = -1040.6
$ws.Range("H55").Value = 2289.2
$ws.Range("I55").Value = 1450
$ws.Range("J55").Value = 2499
$ws.Range("K55").Value = 1450
$ws.Range("L55").Value = 2499
$ws.Range("M55").Value = -1277
$ws.Range("N55").Value = -2845
$ws.Range("H132").Value = 724536.9399999999
$ws.Range("I132").Value = 826790.9
$ws.Range("J132").Value = 8759
$ws.Range("K132").Value = 2480372.7
$ws.Range("L132").Value = 26277
$ws.Range("M132").Value = -2477842.7
$ws.Range("N132").Value = -31337

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = $null
$ws.Range("H81").Value = 3087.889
$ws.Range("I81").Value = 1756.7142
$ws.Range("J81").Value = 7747
$ws.Range("K81").Value = 3513.4284
$ws.Range("L81").Value = 15494
$ws.Range("M81").Value = -2452.4284
$ws.Range("N81").Value = -17616
$ws.Range("H84").Value = 3087.889
$ws.Range("I84").Value = 1756.7142
$ws.Range("J84").Value = 7747
$ws.Range("K84").Value = 17567.142
$ws.Range("L84").Value = 77470
$ws.Range("M84").Value = -12263.142
$ws.Range("N84").Value = -88078
$ws.Range("H122").Value = 2598.05
$ws.Range("I122").Value = 2094.4707
$ws.Range("J122").Value = 5451.6665
$ws.Range("K122").Value = 6283.4121
$ws.Range("L122").Value = 16354.9995
$ws.Range("M122").Value = -3833.4121
$ws.Range("N122").Value = -21254.9995
$ws.Range("H126").Value = 5123.9165
$ws.Range("I126").Value = 6183
$ws.Range("J126").Value = 4064.8333
$ws.Range("K126").Value = 18549
$ws.Range("L126").Value = 12194.4999
$ws.Range("M126").Value = -16079
$ws.Range("N126").Value = -17134.4999
$ws.Range("H132").Value = 587967.7
$ws.Range("I132").Value = 769296.1
$ws.Range("J132").Value = 5126.2144
$ws.Range("K132").Value = 2307888.3
$ws.Range("L132").Value = 15378.6432
$ws.Range("M132").Value = -2305358.3
$ws.Range("N132").Value = -20438.6432
